# "fixed related ground generation"
#
# On the ADC_100MS sheet, a new "clk" pin row is inserted right after the
# existing pin table (between row 18 "v_in2" and the old blank separator
# row), which pushes the separator row and the instance table down by one
# row. The instance table's own content (mem_i0/mem_i1/mem_i2) is unchanged,
# it just slides from rows 22-25 to rows 23-26 as a consequence of the
# insert. Selection moves to C19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADC_100MS")

# Insert a new row at row 19 - shifts the old row 19 (blank separator) and
# everything below it (the instance table) down by one row, matching the
# diff's row 19->20 and 22-25->23-26 shifts.
$ws.Rows.Item(19).EntireRow.Insert()

# The freshly inserted row 19 has no formatting yet; pick it up from row 18
# (the last row of the pin table), which carries the styles the new clk
# row should use (s7/s0/s0/s0/s1/s1/s8).
$ws.Range("B18:H18").Copy()
$ws.Range("B19:H19").PasteSpecial(-4122)

# Fill in the new "clk" pin entry.
$ws.Range("B19").Value = "clk"
$ws.Range("C19").Value = "clk_in"
$ws.Range("D19").Value = "gndd"
$ws.Range("E19").Value = "vddd"
$ws.Range("F19").Formula = "=F18+10"
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = "digital clock"

# Match the author's new selection.
$ws.Range("C19").Select()
